$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 2019 data row (row 31) for OOSS y huelgas annual report
$ws.Range("A31").Value = 2019
$ws.Range("B31").Value = 11926
$ws.Range("C31").Value = 1193104
$ws.Range("D31").Value = 9107664.122137405
$ws.Range("E31").Value = 13.1
$ws.Range("F31").Value = 7503798.742138364
$ws.Range("G31").Value = 15.9
$ws.Range("H31").Value = 1066996
$ws.Range("I31").Value = 5645481.481481481
$ws.Range("J31").Value = 18.9
